$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - index 1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 26433
$ws1.Range("F6").Value = 597
$ws1.Range("F7").Value = 176
$ws1.Range("F10").Value = 354
$ws1.Range("F13").Value = 48
$ws1.Range("F14").Value = 295
$ws1.Range("F15").Value = 51
$ws1.Range("F16").Value = 384
$ws1.Range("F18").Value = 1516
$ws1.Range("F19").Value = 189

# Sheet "演出" (Performance) - index 2
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 4501
$ws2.Range("F6").Value = 188

# Sheet "本地生活" (Local Life) - index 3
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 5016

# Sheet "全部类型" (All Types) - index 4
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 5016
$ws4.Range("F5").Value = 26433
$ws4.Range("F7").Value = 4501
$ws4.Range("F10").Value = 597
$ws4.Range("F13").Value = 176
$ws4.Range("F14").Value = 188
$ws4.Range("F15").Value = 188
$ws4.Range("F23").Value = 354
$ws4.Range("F26").Value = 48
$ws4.Range("F28").Value = 295
$ws4.Range("F29").Value = 51
$ws4.Range("F32").Value = 384
$ws4.Range("F35").Value = 1516
$ws4.Range("F36").Value = 189
$ws4.Range("F38").Value = 34
